$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.292.78'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.671.22'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '682.47'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.49'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.09%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.19%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.97'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.98%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.65%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.290.52'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.04'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.685.96'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.257.30'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '15.75'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.58%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -4.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '470.68'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.90'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.97%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '79.85'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.818.71'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -5.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.86'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.64%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -5.48%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.39%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.48%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.80'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.49'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.96%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -7.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.649.35'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.156'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.12'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.36%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0893'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.37%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.939'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '165.56'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.56'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.000277'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.35%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.09'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.74'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.71'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.86%  '
